# The post "「クリエイティブになるには？」" (row 406) was removed from the
# spreadsheet. Deleting the entire row shifts every row below it up by one,
# which matches the rest of the diff (all subsequent rows renumbered -1 and
# the sheet's used-range dimension shrinking from C576 to C575).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(406).Delete()
